$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add two new header cells (D1, E1), same style as C1 ---
$ws.Range("D1").Value = "SO CAU DUNG"
$ws.Range("E1").Value = "TONG CAU"
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data row: replace MSSV / exam-code text values ---
# A2/B2 must stay text (like the original "752100"/"106"), not auto-convert
# to numbers, so force Text format first, write, then drop back to the
# workbook's default "Normal" style (keeps the values textual without
# leaving a bold/border style behind on the cells).
$ws.Range("A2:B2").NumberFormat = "@"
$ws.Range("A2").Value = "301111"
$ws.Range("B2").Value = "110"
$ws.Range("A2:B2").Style = "Normal"

# --- Data row: numeric score columns ---
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 50
$ws.Range("E2").Value = 50
